$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.656.87"
$ws.Range("E2").Value = "  -1.59%  "

$ws.Range("D3").Value = "1.882.28"
$ws.Range("E3").Value = "  -1.17%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "331.02"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.10%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.08%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4713"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  +1.80%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3964"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.26%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "48.47"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -6.31%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.08066"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.22%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.027"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.50%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "21.82"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +0.07%  "

$ws.Range("D13").Value = "1.882.53"
$ws.Range("E13").Value = "  -1.11%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "5.963"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.94%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.196"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -2.29%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.003"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "86.89"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -3.00%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.00001043"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -2.08%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.06600"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.20"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -3.25%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").Value = "27.666.34"
$ws.Range("E22").Value = "  -1.47%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.513"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.94%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "10.99"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.305"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "

$ws.Range("D26").Value = "2.103.23"
$ws.Range("E26").Value = "  -1.15%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "154.96"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +0.22%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "20.25"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +1.11%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.098"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.25%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.588"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.18%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "122.64"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -2.01%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.09532"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.9622"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.68%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.476"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +0.60%  "

$ws.Range("E35").Value = "  -0.27%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "5.298"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -3.60%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.06119"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -0.19%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02256"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.69%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.225"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.83%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "8.178"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -6.01%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.6009"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -1.92%  "

$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.002"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.1899"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.10%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "10.33"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -4.65%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.5712"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.38%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.248"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -4.11%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "12.14"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -5.00%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "3.412"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.57%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.938"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -3.50%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.06826"
$cell.Style = "Normal"

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "110.27"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "
